$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Baptized-Email"), shifting
# all existing columns C..P to D..Q.
$ws.Columns("C:C").Insert()

# New header for the inserted column.
$ws.Range("C1").Value = "Baptized-Email"

# Row 20 (Winn19 / Tom19) data tweaks:
#  - the "Baptism-Date" value moves from 43565 -> 43527
#  - the Godmother-Email hyperlink text changes to a real address
$ws.Range("E20").Value = 43527
$ws.Range("H20").Value = "alfmat@live.unc.edu"

# Rebuild the hyperlinks: column insert does not shift existing
# hyperlink ranges, so drop them all and recreate at the right spots.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:mattsebastian@live.com")
$ws.Hyperlinks.Add($ws.Range("H20"), "mailto:mattsebastian@live.com")

# Adding a hyperlink re-styles the cell with a freshly minted style;
# restore the original shared hyperlink-cell formatting (style used by
# every other cell in the column) by copying it over.
$ws.Range("H3").Copy()
$ws.Range("H2").PasteSpecial(-4122)
$ws.Range("H3").Copy()
$ws.Range("H20").PasteSpecial(-4122)

# Update the active selection.
$ws.Range("H20").Select()
